$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-06 19:17:43'
$ws.Range("H2").Value = '86%'
$ws.Range("E3").Value = '2026-02-06 19:17:45'
$ws.Range("E4").Value = '2026-02-06 19:17:48'
$ws.Range("J4").Value = '997.4 hPa'
$ws.Range("E5").Value = '2026-02-06 19:17:50'
$ws.Range("J5").Value = '997.6 hPa'
$ws.Range("E6").Value = '2026-02-06 19:17:53'
$ws.Range("J6").Value = '998.8 hPa'
$ws.Range("O6").Value = '15.4 °C'
$ws.Range("E7").Value = '2026-02-06 19:17:55'
$ws.Range("J7").Value = '998.4 hPa'
$ws.Range("E8").Value = '2026-02-06 19:17:58'
$ws.Range("K8").Value = '11.8 MJ/m2'
$ws.Range("O8").Value = '10.2 °C'
$ws.Range("E9").Value = '2026-02-06 19:18:00'
$ws.Range("O9").Value = '5.0 °C'
$ws.Range("E10").Value = '2026-02-06 19:18:03'
$ws.Range("O10").Value = '9.5 °C'
$ws.Range("E11").Value = '2026-02-06 19:18:05'
$ws.Range("J11").Value = '998.8 hPa'
$ws.Range("E12").Value = '2026-02-06 19:18:08'
$ws.Range("H12").Value = '60%'
$ws.Range("N12").Value = '7.7 °C 18:51 TU'
$ws.Range("O12").Value = '14.0 °C'
$ws.Range("E13").Value = '2026-02-06 19:18:10'
$ws.Range("E14").Value = '2026-02-06 19:18:13'
$ws.Range("O14").Value = '-4.2 °C'
$ws.Range("E15").Value = '2026-02-06 19:18:15'
$ws.Range("H15").Value = '71%'
$ws.Range("J15").Value = '997.8 hPa'
$ws.Range("E16").Value = '2026-02-06 19:18:18'
$ws.Range("H16").Value = '83%'
$ws.Range("E17").Value = '2026-02-06 19:18:20'
$ws.Range("J17").Value = '999.0 hPa'
$ws.Range("E18").Value = '2026-02-06 19:18:23'
$ws.Range("I18").Value = '0.2 mm'
$ws.Range("E19").Value = '2026-02-06 19:18:25'
$ws.Range("J19").Value = '999.9 hPa'
$ws.Range("E20").Value = '2026-02-06 19:18:28'
$ws.Range("O20").Value = '-2.0 °C'
$ws.Range("E21").Value = '2026-02-06 19:18:30'
$ws.Range("J21").Value = '998.0 hPa'
$ws.Range("O21").Value = '8.7 °C'
$ws.Range("E22").Value = '2026-02-06 19:18:33'
$ws.Range("H22").Value = '78%'
$ws.Range("E23").Value = '2026-02-06 19:18:35'
$ws.Range("J23").Value = '997.8 hPa'
$ws.Range("E24").Value = '2026-02-06 19:18:38'
$ws.Range("H24").Value = '65%'
$ws.Range("J24").Value = '997.2 hPa'
$ws.Range("O24").Value = '13.0 °C'
$ws.Range("E25").Value = '2026-02-06 19:18:40'
$ws.Range("H25").Value = '76%'
$ws.Range("J25").Value = '998.5 hPa'
$ws.Range("E26").Value = '2026-02-06 19:18:43'
$ws.Range("E27").Value = '2026-02-06 19:18:45'
$ws.Range("H27").Value = '78%'
$ws.Range("J27").Value = '997.8 hPa'
$ws.Range("E28").Value = '2026-02-06 19:18:48'
$ws.Range("J28").Value = '999.9 hPa'
$ws.Range("E29").Value = '2026-02-06 19:18:50'
$ws.Range("H29").Value = '61%'
$ws.Range("O29").Value = '12.7 °C'
$ws.Range("E30").Value = '2026-02-06 19:18:53'
$ws.Range("H30").Value = '77%'
$ws.Range("K30").Value = '9.0 MJ/m2'
$ws.Range("E31").Value = '2026-02-06 19:18:55'
$ws.Range("J31").Value = '999.2 hPa'
$ws.Range("E32").Value = '2026-02-06 19:18:57'
$ws.Range("J32").Value = '999.1 hPa'
$ws.Range("E33").Value = '2026-02-06 19:19:00'
$ws.Range("E34").Value = '2026-02-06 19:19:02'
$ws.Range("E35").Value = '2026-02-06 19:19:05'
$ws.Range("G35").Value = '201 cm'
$ws.Range("E36").Value = '2026-02-06 19:19:07'
$ws.Range("H36").Value = '64%'
$ws.Range("I36").Value = '1.6 mm'
$ws.Range("J36").Value = '1000.0 hPa'
